# Generate Report for Handoff
# - Update localization status from "In Translation" to "Ready for handoff"
# - Refresh the associated handoff timestamps
# - Widen the Status / Latest Handoff Datetime columns to fit the new text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status text: "In Translation" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# Updated handoff timestamps
$wsOverview.Range("G2").Value = "2016-08-24 02:39:03"
$wsDeDe.Range("H2").Value     = "2016-08-24 02:39:03"
$wsZhCn.Range("H2").Value     = "2016-08-24 02:38:56"

# Widen the columns that now hold the longer "Ready for handoff" status text
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.33
